$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ D = 38.82603207464188;  E = -1188835.011349667 }
    3  = @{ D = 33.82565791732615;  E = -1014541.525934224 }
    4  = @{ D = 30.99823694430878;  E = -614977.330991718 }
    5  = @{ D = 16.22681908322957;  E = -309615.2685212554 }
    6  = @{ D = 14.13922243651796;  E = -219310.2006919099 }
    7  = @{ D = 35.87151118214025;  E = -947682.1563098977 }
    8  = @{ D = 30.69634077797294;  E = -809515.3688860169 }
    9  = @{ D = 29.45350707573782;  E = -411545.6223935362 }
    10 = @{ D = 16.75092434359108;  E = -256658.0964985974 }
    11 = @{ D = 14.54501201713191;  E = -188740.3607194683 }
    12 = @{ D = 78.39259224291204;  E = -1975429.738784333 }
    13 = @{ D = 74.03363892612992;  E = -1951102.050533184 }
    14 = @{ D = 52.77572075516835;  E = -646586.4599556178 }
    15 = @{ D = 50.00949033119835;  E = -687173.1145867003 }
    16 = @{ D = 12.83670887836446;  E = -285754.0364846688 }
    17 = @{ D = 109.3941818438427;  E = -1667045.262704854 }
    18 = @{ D = 98.69350734799384;  E = -1650101.137704854 }
    19 = @{ D = 71.24359639498253;  E = -1395774.11472019 }
    20 = @{ D = 46.45385484777381;  E = -544706.7459205341 }
    21 = @{ D = 16.2245674276828;   E = -219015.1786025419 }
}

foreach ($row in $data.Keys) {
    $ws.Range("D$row").Value = $data[$row].D
    $ws.Range("E$row").Value = $data[$row].E
}
